$d = $word.ActiveDocument

# --- Change 1 ---
# "...framework of hierarchical structure formation. That is, that most..."
# becomes
# "...framework of hierarchical structure formation, which occurs due to
#  the assumption of gravitational instability drawing structures together.
#  That is, that most..."
$d.Content.Find.Execute(
    "hierarchical structure formation. That is,", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "hierarchical structure formation, which occurs due to the assumption of gravitational instability drawing structures together. That is,",
    2) | Out-Null

# --- Change 2 ---
# "The push towards investigating this is the 'spin crisis' ..."
# becomes
# "The push towards investigating the process by which angular momentum is
#  formed is the 'spin crisis' ..."
$d.Content.Find.Execute(
    "investigating this is the", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "investigating the process by which angular momentum is formed is the",
    2) | Out-Null

# --- Change 3 / 4 ---
# Move the (hidden) "_GoBack" bookmark from the empty paragraph right
# before "METHOD:" into the "(CORRELATION WITH COSMIC WEB)" paragraph,
# splitting its text right after "CORRELATI" -> "(CORRELATI" | "ON WITH COSMIC WEB)"
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Find the "(CORRELATION WITH COSMIC WEB)" paragraph and compute the split point.
$r = $d.Content
$r.Find.Execute("(CORRELATION WITH COSMIC WEB)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitPos = $r.Start + "(CORRELATI".Length
$bookmarkRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange) | Out-Null
